$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.262.18"
$ws.Range("E2").Value = "'  -1.23%  "
$ws.Range("D3").Value = "'2.249.06"
$ws.Range("E3").Value = "'  -1.23%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'247.47"
$ws.Range("E5").Value = "'  -1.48%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = "'  -2.39%  "
$ws.Range("D7").Value = "'74.76"
$ws.Range("E7").Value = "'  -1.33%  "
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E9").Value = "'  -3.75%  "
$ws.Range("D10").Value = "'42.37"
$ws.Range("E10").Value = "'  +6.12%  "
$ws.Range("D11").Value = "'0.0939"
$ws.Range("E11").Value = "'  -3.68%  "
$ws.Range("E12").Value = "'  -2.64%  "
$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "'  -2.55%  "
$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'2.583.05"
$ws.Range("E14").Value = "'  -1.45%  "
$ws.Range("B15").Value = "'Chainlink"
$ws.Range("C15").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'14.50"
$ws.Range("E15").Value = "'  -3.93%  "
$ws.Range("B16").Value = "'Polygon"
$ws.Range("C16").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.858"
$ws.Range("E16").Value = "'  -1.15%  "
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.248.79"
$ws.Range("E17").Value = "'  -1.65%  "
$ws.Range("B18").Value = "'WrappedBTC"
$ws.Range("C18").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'42.139.81"
$ws.Range("E18").Value = "'  -1.27%  "
$ws.Range("B19").Value = "'ShibaInu"
$ws.Range("C19").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.0₃0983"
$ws.Range("E19").Value = "'  -1.30%  "
$ws.Range("B20").Value = "'Uniswap"
$ws.Range("C20").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.14"
$ws.Range("E20").Value = "'  -1.07%  "
$ws.Range("B21").Value = "'Litecoin"
$ws.Range("C21").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'72.01"
$ws.Range("E21").Value = "'  -0.42%  "
$ws.Range("B22").Value = "'ImmutableX"
$ws.Range("C22").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'2.25"
$ws.Range("E22").Value = "'  +4.41%  "
$ws.Range("B23").Value = "'BitcoinCash"
$ws.Range("C23").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'230.42"
$ws.Range("E23").Value = "'  -2.55%  "
$ws.Range("B24").Value = "'InternetComputer(DFINITY)"
$ws.Range("C24").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "'8.96"
$ws.Range("E24").Value = "'  +39.33%  "
$ws.Range("B25").Value = "'Dai"
$ws.Range("C25").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "'  +0.06%  "
$ws.Range("B26").Value = "'Cosmos"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'11.23"
$ws.Range("E26").Value = "'  -0.14%  "
$ws.Range("B27").Value = "'WEMIXToken"
$ws.Range("C27").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "'3.64"
$ws.Range("E27").Value = "'  -5.65%  "
$ws.Range("B28").Value = "'PancakeSwap"
$ws.Range("C28").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.31"
$ws.Range("E28").Value = "'  -3.27%  "
$ws.Range("B29").Value = "'Monero"
$ws.Range("C29").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'169.58"
$ws.Range("E29").Value = "'  +1.21%  "
$ws.Range("B30").Value = "'Toncoin"
$ws.Range("C30").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.10"
$ws.Range("E30").Value = "'  -1.72%  "
$ws.Range("B31").Value = "'EthereumClassic"
$ws.Range("C31").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'20.72"
$ws.Range("E31").Value = "'  -1.25%  "
$ws.Range("B32").Value = "'Hedera"
$ws.Range("C32").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0828"
$ws.Range("E32").Value = "'  -3.74%  "
$ws.Range("B33").Value = "'Kaspa"
$ws.Range("C33").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.119"
$ws.Range("E33").Value = "'  -5.19%  "
$ws.Range("B34").Value = "'InjectiveProtocol"
$ws.Range("C34").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'30.40"
$ws.Range("E34").Value = "'  -4.96%  "
$ws.Range("B35").Value = "'Stellar"
$ws.Range("C35").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.125"
$ws.Range("E35").Value = "'  -2.06%  "
$ws.Range("B36").Value = "'Filecoin"
$ws.Range("C36").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'5.22"
$ws.Range("E36").Value = "'  +9.85%  "
$ws.Range("B37").Value = "'RenderToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.55"
$ws.Range("E37").Value = "'  +0.18%  "
$ws.Range("B38").Value = "'VeChain"
$ws.Range("C38").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0305"
$ws.Range("E38").Value = "'  -0.27%  "
$ws.Range("B39").Value = "'Celestia"
$ws.Range("C39").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'13.51"
$ws.Range("E39").Value = "'  +0.00%  "
$ws.Range("B40").Value = "'LidoDAOToken"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.19"
$ws.Range("E40").Value = "'  -4.00%  "
$ws.Range("B41").Value = "'THORChain"
$ws.Range("C41").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'5.81"
$ws.Range("E41").Value = "'  -1.27%  "
$ws.Range("B42").Value = "'Aave"
$ws.Range("C42").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'109.65"
$ws.Range("E42").Value = "'  +3.09%  "
$ws.Range("B43").Value = "'MultiversX"
$ws.Range("C43").Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'61.60"
$ws.Range("E43").Value = "'  -0.10%  "
$ws.Range("B44").Value = "'Algorand"
$ws.Range("C44").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.202"
$ws.Range("E44").Value = "'  -2.44%  "
$ws.Range("B45").Value = "'FraxShare"
$ws.Range("C45").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'8.63"
$ws.Range("E45").Value = "'  -3.30%  "
$ws.Range("B46").Value = "'Cronos"
$ws.Range("C46").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.101"
$ws.Range("E46").Value = "'  +0.52%  "
$ws.Range("B47").Value = "'BinanceUSD"
$ws.Range("C47").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").Value = "'0.997"
$ws.Range("E47").Value = "'  -0.36%  "
$ws.Range("B48").Value = "'ARBITRUM"
$ws.Range("C48").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.12"
$ws.Range("E48").Value = "'  -3.82%  "
$ws.Range("E49").Value = "'  -0.52%  "
$ws.Range("B50").Value = "'NEARProtocol"
$ws.Range("C50").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.32"
$ws.Range("E50").Value = "'  +2.93%  "
$ws.Range("B51").Value = "'SynthetixNetwork"
$ws.Range("C51").Value = "'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'4.13"
$ws.Range("E51").Value = "'  -1.81%  "
